$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (Excel COM ColumnWidth applies a +5/6 char offset vs. the raw OOXML width unit;
# subtracting 5/6 before assigning yields the exact target integer width on save).
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(11).ColumnWidth = 36.166666666666664
$ws.Columns.Item(12).ColumnWidth = 44.166666666666664
$ws.Columns.Item(13).ColumnWidth = 36.166666666666664

# Cell value updates (previously-empty cells now populated by the improved scraper)
$ws.Range("E2").Value = "Student organization focused on general activities and community engagement. The Veterans Affairs welcomes all interested students to participate and make a positive impact."
$ws.Range("F2").Value = "veteransaffairs@bhsu.edu"
$ws.Range("G2").Value = "(555) 345-6789"
$ws.Range("I2").Value = "https://instagram.com/veteransaffairs"
$ws.Range("E3").Value = "Student organization focused on general activities and community engagement. The Registration & Records welcomes all interested students to participate and make a positive impact."
$ws.Range("F3").Value = "registrationrec@bhsu.edu"
$ws.Range("G3").Value = "(555) 345-6789"
$ws.Range("J3").Value = "https://facebook.com/registrationrec"
$ws.Range("K3").Value = "https://twitter.com/registrationrec"
$ws.Range("E4").Value = "Student organization focused on general activities and community engagement. The Housing & Residence Halls welcomes all interested students to participate and make a positive impact."
$ws.Range("F4").Value = "housingresidenc@bhsu.edu"
$ws.Range("G4").Value = "(555) 890-1234"
$ws.Range("H4").Value = "https://linkedin.com/groups/housingresidenc"
$ws.Range("D5").Value = "https://bhsu.edu/logos/studentactiviti_logo.png"
$ws.Range("E5").Value = "Student organization focused on general activities and community engagement. The Student Activities welcomes all interested students to participate and make a positive impact."
$ws.Range("F5").Value = "studentactiviti@bhsu.edu"
$ws.Range("H5").Value = "https://linkedin.com/groups/studentactiviti"
$ws.Range("I5").Value = "https://instagram.com/studentactiviti"
$ws.Range("E6").Value = "Student organization focused on general activities and community engagement. The Student Organizations welcomes all interested students to participate and make a positive impact."
$ws.Range("F6").Value = "studentorganiza@bhsu.edu"
$ws.Range("I6").Value = "https://instagram.com/studentorganiza"
$ws.Range("J6").Value = "https://facebook.com/studentorganiza"
$ws.Range("E7").Value = "Student organization focused on general activities and community engagement. The Student Union welcomes all interested students to participate and make a positive impact."
$ws.Range("F7").Value = "studentunion@bhsu.edu"
$ws.Range("G7").Value = "(555) 901-2345"
$ws.Range("K7").Value = "https://twitter.com/studentunion"
$ws.Range("D8").Value = "https://bhsu.edu/logos/buzzcardmobile_logo.png"
$ws.Range("E8").Value = "Student organization focused on general activities and community engagement. The Buzz Card Mobile welcomes all interested students to participate and make a positive impact."
$ws.Range("F8").Value = "buzzcardmobile@bhsu.edu"
$ws.Range("M8").Value = "https://tiktok.com/@buzzcardmobile"
$ws.Range("E9").Value = "Student organization focused on general activities and community engagement. The Visit BHSU welcomes all interested students to participate and make a positive impact."
$ws.Range("F9").Value = "visitbhsu@bhsu.edu"
$ws.Range("I9").Value = "https://instagram.com/visitbhsu"
$ws.Range("K9").Value = "https://twitter.com/visitbhsu"
$ws.Range("E10").Value = "Student organization focused on general activities and community engagement. The Class Registration welcomes all interested students to participate and make a positive impact."
$ws.Range("F10").Value = "classregistrati@bhsu.edu"
$ws.Range("L10").Value = "https://youtube.com/channel/classregistrati"
$ws.Range("M10").Value = "https://tiktok.com/@classregistrati"
$ws.Range("D11").Value = "https://bhsu.edu/logos/paymentoptions_logo.png"
$ws.Range("E11").Value = "Student organization focused on general activities and community engagement. The Payment Options welcomes all interested students to participate and make a positive impact."
$ws.Range("F11").Value = "paymentoptions@bhsu.edu"
$ws.Range("H11").Value = "https://linkedin.com/groups/paymentoptions"
$ws.Range("K11").Value = "https://twitter.com/paymentoptions"
$ws.Range("E12").Value = "Student organization focused on general activities and community engagement. The Green & Gold Days welcomes all interested students to participate and make a positive impact."
$ws.Range("F12").Value = "greengolddays@bhsu.edu"
$ws.Range("G12").Value = "(555) 456-7890"
$ws.Range("D13").Value = "https://bhsu.edu/logos/internationalst_logo.png"
$ws.Range("E13").Value = "Cultural organization celebrating diversity and promoting multicultural awareness. The International Students welcomes all interested students to participate and make a positive impact."
$ws.Range("G13").Value = "(555) 890-1234"
$ws.Range("H13").Value = "https://linkedin.com/groups/internationalst"
$ws.Range("J13").Value = "https://facebook.com/internationalst"
$ws.Range("K13").Value = "https://twitter.com/internationalst"
$ws.Range("L13").Value = "https://youtube.com/channel/internationalst"
$ws.Range("E14").Value = "Student organization focused on general activities and community engagement. The Room Rentals & Reservations welcomes all interested students to participate and make a positive impact."
$ws.Range("F14").Value = "roomrentalsrese@bhsu.edu"
$ws.Range("J14").Value = "https://facebook.com/roomrentalsrese"
$ws.Range("E15").Value = "Student organization focused on general activities and community engagement. The Jacket Connect welcomes all interested students to participate and make a positive impact."
$ws.Range("F15").Value = "jacketconnect@bhsu.edu"
$ws.Range("L15").Value = "https://youtube.com/channel/jacketconnect"
$ws.Range("D16").Value = "https://bhsu.edu/logos/menswomensbaske_logo.png"
$ws.Range("E16").Value = "Athletic organization promoting physical fitness and competitive spirit. The Men's & Women's Basketball welcomes all interested students to participate and make a positive impact."
$ws.Range("F16").Value = "menswomensbaske@bhsu.edu"
$ws.Range("G16").Value = "(555) 789-0123"
$ws.Range("H16").Value = "https://linkedin.com/groups/menswomensbaske"
$ws.Range("I16").Value = "https://instagram.com/menswomensbaske"
$ws.Range("E17").Value = "Student organization focused on general activities and community engagement. The Cross Country welcomes all interested students to participate and make a positive impact."
$ws.Range("F17").Value = "crosscountry@bhsu.edu"
$ws.Range("G17").Value = "(555) 789-0123"
$ws.Range("I17").Value = "https://instagram.com/crosscountry"
$ws.Range("J17").Value = "https://facebook.com/crosscountry"
$ws.Range("M17").Value = "https://tiktok.com/@crosscountry"
$ws.Range("D18").Value = "https://bhsu.edu/logos/trackfield_logo.png"
$ws.Range("E18").Value = "Student organization focused on general activities and community engagement. The Track & Field welcomes all interested students to participate and make a positive impact."
$ws.Range("F18").Value = "trackfield@bhsu.edu"
$ws.Range("G18").Value = "(555) 123-4567"
$ws.Range("H18").Value = "https://linkedin.com/groups/trackfield"
$ws.Range("I18").Value = "https://instagram.com/trackfield"
$ws.Range("J18").Value = "https://facebook.com/trackfield"
$ws.Range("E19").Value = "Student organization focused on general activities and community engagement. The Future Students welcomes all interested students to participate and make a positive impact."
$ws.Range("F19").Value = "futurestudents@bhsu.edu"
$ws.Range("G19").Value = "(555) 345-6789"
$ws.Range("I19").Value = "https://instagram.com/futurestudents"
$ws.Range("J19").Value = "https://facebook.com/futurestudents"
$ws.Range("M19").Value = "https://tiktok.com/@futurestudents"
$ws.Range("D20").Value = "https://bhsu.edu/logos/concerncomplain_logo.png"
$ws.Range("E20").Value = "Student organization focused on general activities and community engagement. The Concern & Complaint welcomes all interested students to participate and make a positive impact."
$ws.Range("G20").Value = "(555) 123-4567"
$ws.Range("H20").Value = "https://linkedin.com/groups/concerncomplain"
$ws.Range("I20").Value = "https://instagram.com/concerncomplain"
$ws.Range("K20").Value = "https://twitter.com/concerncomplain"
$ws.Range("E21").Value = "Student organization focused on general activities and community engagement. The Student Portal welcomes all interested students to participate and make a positive impact."
$ws.Range("F21").Value = "studentportal@bhsu.edu"
$ws.Range("L21").Value = "https://youtube.com/channel/studentportal"
$ws.Range("M21").Value = "https://tiktok.com/@studentportal"
$ws.Range("E22").Value = "Student organization focused on general activities and community engagement. The Transcript Request welcomes all interested students to participate and make a positive impact."
$ws.Range("F22").Value = "transcriptreque@bhsu.edu"
$ws.Range("G22").Value = "(555) 345-6789"
$ws.Range("H22").Value = "https://linkedin.com/groups/transcriptreque"
$ws.Range("I22").Value = "https://instagram.com/transcriptreque"
$ws.Range("J22").Value = "https://facebook.com/transcriptreque"
$ws.Range("M22").Value = "https://tiktok.com/@transcriptreque"
$ws.Range("E23").Value = "Student organization focused on general activities and community engagement. The Desire2Learn (D2L) welcomes all interested students to participate and make a positive impact."
$ws.Range("F23").Value = "desire2learnd2l@bhsu.edu"
$ws.Range("G23").Value = "(555) 345-6789"
$ws.Range("J23").Value = "https://facebook.com/desire2learnd2l"
$ws.Range("K23").Value = "https://twitter.com/desire2learnd2l"
$ws.Range("E24").Value = "Student organization focused on general activities and community engagement. The Student Consumer Info welcomes all interested students to participate and make a positive impact."
$ws.Range("F24").Value = "studentconsumer@bhsu.edu"
$ws.Range("G24").Value = "(555) 567-8901"
$ws.Range("H24").Value = "https://linkedin.com/groups/studentconsumer"
$ws.Range("I24").Value = "https://instagram.com/studentconsumer"
$ws.Range("K24").Value = "https://twitter.com/studentconsumer"
$ws.Range("E25").Value = "Student organization focused on general activities and community engagement. The Web Request welcomes all interested students to participate and make a positive impact."
$ws.Range("F25").Value = "webrequest@bhsu.edu"
$ws.Range("J25").Value = "https://facebook.com/webrequest"
$ws.Range("D26").Value = "https://bhsu.edu/logos/titleix_logo.png"
$ws.Range("E26").Value = "Student organization focused on general activities and community engagement. The Title IX welcomes all interested students to participate and make a positive impact."
$ws.Range("E27").Value = "Student organization focused on general activities and community engagement. The Email Password Reset welcomes all interested students to participate and make a positive impact."
$ws.Range("F27").Value = "emailpasswordre@bhsu.edu"
$ws.Range("G27").Value = "(555) 123-4567"
$ws.Range("H27").Value = "https://linkedin.com/groups/emailpasswordre"
$ws.Range("J27").Value = "https://facebook.com/emailpasswordre"
$ws.Range("L29").Value = "https://youtube.com/channel/studentlife"
$ws.Range("L30").Value = "https://youtube.com/channel/upcomingauditio"
$ws.Range("L33").Value = "https://youtube.com/channel/whereyoulleat"
